$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected error in s_profile
# stratification_type value: "step" -> "uniform"
$ws.Range("B17").Value = "uniform"

# rhoa0 value: 1.05 -> 0.998
$ws.Range("B8").Value = 0.998

# rhoa_upper value: 1.01 -> 0.998
$ws.Range("B9").Value = 0.998

# update active selection to B7
$ws.Range("B7").Select()
